$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("O20").Value = '[''Portugal'', ''Northern Ireland'']'
$ws.Range("H46").Value = '[''Russia'', 3, -4, 2]'
$ws.Range("M46").Value = '[''Switzerland'', ''Portugal'', ''Ukraine'', ''Russia'']'
$ws.Range("P46").Value = 0
$ws.Range("Q46").Value = 4
$ws.Range("H47").Value = '[''Russia'', 3, -4, 2]'
$ws.Range("M47").Value = '[''Switzerland'', ''Portugal'', ''Ukraine'', ''Russia'']'
$ws.Range("Q47").Value = 4
$ws.Range("H48").Value = '[''Russia'', 3, -5, 2]'
$ws.Range("M48").Value = '[''Switzerland'', ''Portugal'', ''Ukraine'', ''Russia'']'
$ws.Range("Q48").Value = 4
$ws.Range("H49").Value = '[''Russia'', 3, -5, 2]'
$ws.Range("M49").Value = '[''Switzerland'', ''Portugal'', ''Ukraine'', ''Russia'']'
$ws.Range("Q49").Value = 4
$ws.Range("H50").Value = '[''Russia'', 3, -5, 2]'
$ws.Range("M50").Value = '[''Switzerland'', ''Portugal'', ''Ukraine'', ''Russia'']'
$ws.Range("Q50").Value = 4
$ws.Range("H51").Value = '[''Russia'', 3, -5, 2]'
$ws.Range("N51").Value = '[''Russia'', ''Spain'']'
$ws.Range("Q51").Value = 5
$ws.Range("H52").Value = '[''Russia'', 3, -5, 2]'
$ws.Range("M52").Value = '[''Switzerland'', ''Portugal'', ''Ukraine'', ''Russia'']'
$ws.Range("Q52").Value = 6
$ws.Range("H53").Value = '[''Russia'', 3, -5, 2]'
$ws.Range("N53").Value = '[''Russia'', ''Spain'']'
$ws.Range("Q53").Value = 7
$ws.Range("H54").Value = '[''Russia'', 3, -5, 2]'
$ws.Range("N54").Value = '[''Russia'', ''Spain'']'
$ws.Range("Q54").Value = 8
$ws.Range("H55").Value = '[''Russia'', 3, -5, 2]'
$ws.Range("N55").Value = '[''Ukraine'', ''Russia'']'
$ws.Range("Q55").Value = 9
$ws.Range("H56").Value = '[''Russia'', 3, -5, 2]'
$ws.Range("N56").Value = '[''Ukraine'', ''Russia'']'
$ws.Range("Q56").Value = 9
$ws.Range("H57").Value = '[''Russia'', 3, -5, 2]'
$ws.Range("N57").Value = '[''Slovakia'', ''Russia'']'
$ws.Range("Q57").Value = 10
$ws.Range("H58").Value = '[''Russia'', 3, -5, 2]'
$ws.Range("N58").Value = '[''Slovakia'', ''Russia'']'
$ws.Range("Q58").Value = 10
$ws.Range("H59").Value = '[''Russia'', 3, -5, 2]'
$ws.Range("N59").Value = '[''Slovakia'', ''Russia'']'
$ws.Range("Q59").Value = 10
$ws.Range("H60").Value = '[''Russia'', 3, -5, 2]'
$ws.Range("N60").Value = '[''Slovakia'', ''Russia'']'
$ws.Range("Q60").Value = 10
$ws.Range("H61").Value = '[''Russia'', 3, -5, 2]'
$ws.Range("N61").Value = '[''Slovakia'', ''Russia'']'
$ws.Range("Q61").Value = 10
$ws.Range("H62").Value = '[''Russia'', 3, -5, 2]'
$ws.Range("N62").Value = '[''Slovakia'', ''Russia'']'
$ws.Range("Q62").Value = 10
$ws.Range("H63").Value = '[''Russia'', 3, -5, 2]'
$ws.Range("N63").Value = '[''Russia'', ''Slovakia'']'
$ws.Range("O63").Value = '[''Russia'', ''Slovakia'']'
$ws.Range("Q63").Value = 10
$ws.Range("H64").Value = '[''Russia'', 3, -5, 2]'
$ws.Range("N64").Value = '[''Russia'', ''Slovakia'']'
$ws.Range("O64").Value = '[''Russia'', ''Slovakia'']'
$ws.Range("Q64").Value = 10
$ws.Range("H65").Value = '[''Russia'', 3, -5, 2]'
$ws.Range("N65").Value = '[''Russia'', ''Slovakia'']'
$ws.Range("O65").Value = '[''Russia'', ''Slovakia'']'
$ws.Range("Q65").Value = 10
$ws.Range("H66").Value = '[''Russia'', 3, -5, 2]'
$ws.Range("N66").Value = '[''Russia'', ''Slovakia'']'
$ws.Range("O66").Value = '[''Russia'', ''Slovakia'']'
$ws.Range("Q66").Value = 10
$ws.Range("H67").Value = '[''Russia'', 3, -5, 2]'
$ws.Range("N67").Value = '[''Russia'', ''Slovakia'']'
$ws.Range("O67").Value = '[''Russia'', ''Slovakia'']'
$ws.Range("Q67").Value = 11
$ws.Range("H68").Value = '[''Russia'', 3, -5, 2]'
$ws.Range("N68").Value = '[''Russia'', ''Slovakia'']'
$ws.Range("O68").Value = '[''Russia'', ''Slovakia'']'
$ws.Range("Q68").Value = 11
$ws.Range("H69").Value = '[''Russia'', 3, -5, 2]'
$ws.Range("N69").Value = '[''Russia'', ''Slovakia'']'
$ws.Range("O69").Value = '[''Russia'', ''Slovakia'']'
$ws.Range("Q69").Value = 11
$ws.Range("H70").Value = '[''Russia'', 3, -5, 2]'
$ws.Range("N70").Value = '[''Russia'', ''Slovakia'']'
$ws.Range("O70").Value = '[''Russia'', ''Slovakia'']'
$ws.Range("Q70").Value = 12
$ws.Range("H71").Value = '[''Russia'', 3, -5, 2]'
$ws.Range("N71").Value = '[''Russia'', ''Slovakia'']'
$ws.Range("O71").Value = '[''Russia'', ''Slovakia'']'
$ws.Range("Q71").Value = 13
$ws.Range("H72").Value = '[''Russia'', 3, -5, 2]'
$ws.Range("N72").Value = '[''Russia'', ''Slovakia'']'
$ws.Range("O72").Value = '[''Russia'', ''Slovakia'']'
$ws.Range("Q72").Value = 14
$ws.Range("H73").Value = '[''Russia'', 3, -5, 2]'
$ws.Range("N73").Value = '[''Russia'', ''Slovakia'']'
$ws.Range("O73").Value = '[''Russia'', ''Slovakia'']'
$ws.Range("Q73").Value = 15
$ws.Range("H74").Value = '[''Russia'', 3, -5, 2]'
$ws.Range("N74").Value = '[''Russia'', ''Slovakia'']'
$ws.Range("O74").Value = '[''Russia'', ''Slovakia'']'
$ws.Range("Q74").Value = 16
$ws.Range("O98").Value = '[''Netherlands'', ''Georgia'']'
$ws.Range("O99").Value = '[''Netherlands'', ''Georgia'']'
$ws.Range("O100").Value = '[''Netherlands'', ''Georgia'']'
